$d = $word.ActiveDocument
$anchor = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$anchor.Collapse(0)

# Phase 1: create all new (empty) paragraphs, collecting references
$newParas = @()
$anchor.InsertParagraphAfter()
$pObj = $d.Paragraphs.Item($d.Paragraphs.Count)
$newParas += $pObj
$anchor = $pObj.Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$pObj = $d.Paragraphs.Item($d.Paragraphs.Count)
$newParas += $pObj
$anchor = $pObj.Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$pObj = $d.Paragraphs.Item($d.Paragraphs.Count)
$newParas += $pObj
$anchor = $pObj.Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$pObj = $d.Paragraphs.Item($d.Paragraphs.Count)
$newParas += $pObj
$anchor = $pObj.Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$pObj = $d.Paragraphs.Item($d.Paragraphs.Count)
$newParas += $pObj
$anchor = $pObj.Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$pObj = $d.Paragraphs.Item($d.Paragraphs.Count)
$newParas += $pObj
$anchor = $pObj.Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$pObj = $d.Paragraphs.Item($d.Paragraphs.Count)
$newParas += $pObj
$anchor = $pObj.Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$pObj = $d.Paragraphs.Item($d.Paragraphs.Count)
$newParas += $pObj
$anchor = $pObj.Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$pObj = $d.Paragraphs.Item($d.Paragraphs.Count)
$newParas += $pObj
$anchor = $pObj.Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$pObj = $d.Paragraphs.Item($d.Paragraphs.Count)
$newParas += $pObj
$anchor = $pObj.Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$pObj = $d.Paragraphs.Item($d.Paragraphs.Count)
$newParas += $pObj
$anchor = $pObj.Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$pObj = $d.Paragraphs.Item($d.Paragraphs.Count)
$newParas += $pObj
$anchor = $pObj.Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$pObj = $d.Paragraphs.Item($d.Paragraphs.Count)
$newParas += $pObj
$anchor = $pObj.Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$pObj = $d.Paragraphs.Item($d.Paragraphs.Count)
$newParas += $pObj
$anchor = $pObj.Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$pObj = $d.Paragraphs.Item($d.Paragraphs.Count)
$newParas += $pObj
$anchor = $pObj.Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$pObj = $d.Paragraphs.Item($d.Paragraphs.Count)
$newParas += $pObj
$anchor = $pObj.Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$pObj = $d.Paragraphs.Item($d.Paragraphs.Count)
$newParas += $pObj
$anchor = $pObj.Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$pObj = $d.Paragraphs.Item($d.Paragraphs.Count)
$newParas += $pObj
$anchor = $pObj.Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$pObj = $d.Paragraphs.Item($d.Paragraphs.Count)
$newParas += $pObj
$anchor = $pObj.Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$pObj = $d.Paragraphs.Item($d.Paragraphs.Count)
$newParas += $pObj
$anchor = $pObj.Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$pObj = $d.Paragraphs.Item($d.Paragraphs.Count)
$newParas += $pObj
$anchor = $pObj.Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$pObj = $d.Paragraphs.Item($d.Paragraphs.Count)
$newParas += $pObj
$anchor = $pObj.Range
$anchor.Collapse(0)

# Phase 2: fill in text (plain concatenation per paragraph, no bold yet)
$newParas[0].Range.Text = 'Started a user sign in feature. '
$newParas[1].Range.Text = 'Within the let statement I said ‘create :user’ hoping that factory girl would work her magic. In the tutorial it asks to write User.create with all the properties in brackets. May not work! '
$newParas[2].Range.Text = 'Replaced the above with the let statement, no change to error message. '
$newParas[3].Range.Text = ' Error message says “Undefined method ‘name’. I now added some fill_in lines to the sign_in method. '
$newParas[4].Range.Text = 'On track, now says ‘Unable to find field: :username’. This is directing us to write a new route in our controller with an erb file containing the form to sign in. '
$newParas[5].Range.Text = 'Erb file written with form details on. Stored in views/sessions/new'
$newParas[6].Range.Text = 'Now writing first unit test (Not sure why now). Instead of following the tutorial, I’ve changed the let statement so it calls factory girl to create the user. Again, not sure if this works. '
$newParas[7].Range.Text = 'Placed authenticate class method in User model. Required spec helper in new user feature test file. '
$newParas[8].Range.Text = 'First authenticate not passing. Still being instructed to write does not authenticate test. '
$newParas[9].Range.Text = 'Reached a stage where spec/models/user tests file was not running. Didn’t name it user_spec.rb '
$newParas[10].Range.Text = 'Now have two feature tests. One passing, one not. ‘User does not authenticate when given incorrect password’ is not passing. '
$newParas[11].Range.Text = 'The above unit test now passes as we add a conditional specifying if the user exists and Bcrypt stuff. '
$newParas[12].Range.Text = 'User sign in with correct credentials not passing. Trying to figure out why. '
$newParas[13].Range.Text = 'Will take out my factory girl create in the let in the user_management spec and replace it with what the tutorial specifies. '
$newParas[14].Range.Text = 'That alone made no difference to the error message. '
$newParas[15].Range.Text = 'The error message changed. It had undefined method for name, when I should have said user.username, not user.name as I did before. That’s changed, new error message saying “expected to find text “Welcome Adrian1707” in the page. It’s not showing so investigating why. Will rackup and see what happens. '
$newParas[16].Range.Text = 'DAMN. In post ‘/sessions’ it had redirect to ‘/links’ which of course doesn’t exist because we’re doing ‘peeps’ now. That’s corrected. Test still doesn’t pass but getting there. Now at least in the browser it has what we want. Test not passing for silly reason.'
$newParas[17].Range.Text = 'My let statement is not the same as the tutorial in the user spec file. I wonder if that’s a reason, I heavily doubt it though. Changing now. '
$newParas[18].Range.Text = 'Made another 2 tests fail. I guess I’ll leave it for now'
$newParas[19].Range.Text = 'DONE IT! Simply changed the let statement in the user management spec file to create a user from factory girl. They should be equivalent though so I’m not sure why it made a difference. Still confused by this. '
$newParas[20].Range.Text = 'Also still confused by how the self_authenticate method works with BCrypt. '
$newParas[21].Range.Text = 'Moving onto signing out. '

# Phase 3: apply bold to specific sub-ranges, using absolute offsets
$pRange = $newParas[1].Range
$boldRange = $d.Range($pRange.Start + 176, $pRange.Start + 190)
$boldRange.Font.Bold = 1